$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the creditcard.csv row (row 10) entirely
$ws.Rows.Item(10).Delete()

# 2. Rename existing iforestASD headers (G1:J1) to MILOF headers
$ws.Range("G1").Value = "MILOF_identified"
$ws.Range("H1").Value = "MILOF_Overlap_merlin"
$ws.Range("I1").Value = "MILOFbest_param"
$ws.Range("J1").Value = "MILOFtime_taken"

# 3. Add new ARIMAFD headers in K1:N1 and copy header formatting from G1:J1
$ws.Range("K1").Value = "ARIMAFD_identified"
$ws.Range("L1").Value = "ARIMAFD_Overlap_merlin"
$ws.Range("M1").Value = "ARIMAFDbest_param"
$ws.Range("N1").Value = "ARIMAFDtime_taken"
$ws.Range("G1:J1").Copy()
$ws.Range("K1:N1").PasteSpecial(-4122)

# 4. Update MILOF (G,H,I,J) results and add ARIMAFD (K,L,M,N) results for each data row
$ws.Range("G2").Value = "[88, 552, 797, 825]"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "{'Numk': 14, 'KPar': 11, 'Bucket_index': 500}"
$ws.Range("J2").Value = 9.049631076006335
$ws.Range("K2").Value = "[133, 464, 628, 963]"
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = "{'window_size': 92}"
$ws.Range("N2").Value = 1188.046888874

$ws.Range("G3").Value = "[135, 158, 214, 253]"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "{'Numk': 10, 'KPar': 13, 'Bucket_index': 500}"
$ws.Range("J3").Value = 11.40097611300007
$ws.Range("K3").Value = "[314, 325, 684, 968]"
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = "{'window_size': 181}"
$ws.Range("N3").Value = 5028.996044122003

$ws.Range("G4").Value = "[202, 227, 295, 372]"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = "{'Numk': 9, 'KPar': 11, 'Bucket_index': 500}"
$ws.Range("J4").Value = 10.41810625100334
$ws.Range("K4").Value = "[181, 190, 583, 597]"
$ws.Range("L4").Value = 0.25
$ws.Range("M4").Value = "{'window_size': 123}"
$ws.Range("N4").Value = 265.9657750680053

$ws.Range("G5").Value = "[175, 312]"
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = "{'Numk': 11, 'KPar': 5, 'Bucket_index': 500}"
$ws.Range("J5").Value = 6.335740316993906
$ws.Range("K5").Value = "[503, 745]"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "{'window_size': 176}"
$ws.Range("N5").Value = 3866.292378414015

$ws.Range("G6").Value = "[372, 383]"
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = "{'Numk': 10, 'KPar': 3, 'Bucket_index': 500}"
$ws.Range("J6").Value = 5.555181247997098
$ws.Range("K6").Value = "[203, 987]"
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = "{'window_size': 203}"
$ws.Range("N6").Value = 4814.967263817991

$ws.Range("G7").Value = "[58, 64, 528, 536]"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "{'Numk': 6, 'KPar': 6, 'Bucket_index': 500}"
$ws.Range("J7").Value = 5.160258635005448
$ws.Range("K7").Value = "[104, 155, 194, 199]"
$ws.Range("L7").Value = 0.25
$ws.Range("M7").Value = "{'window_size': 90}"
$ws.Range("N7").Value = 962.8695973639842

$ws.Range("G8").Value = "[177, 787, 958, 1068]"
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "{'Numk': 14, 'KPar': 5, 'Bucket_index': 500}"
$ws.Range("J8").Value = 5.478524310994544
$ws.Range("K8").Value = "[435, 444, 571, 628]"
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = "{'window_size': 268}"
$ws.Range("N8").Value = 8245.638174999011

$ws.Range("G9").Value = "[443, 798, 806, 810, 823]"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = "{'Numk': 12, 'KPar': 11, 'Bucket_index': 500}"
$ws.Range("J9").Value = 10.42844949501159
$ws.Range("K9").Value = "[148, 178, 787, 980]"
$ws.Range("L9").Value = 0.25
$ws.Range("M9").Value = "{'window_size': 148}"
$ws.Range("N9").Value = 1719.089206571982

